# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text (prevents Excel from
# auto-converting numeric-looking strings, e.g. "26.854.44" or "3.100",
# into numbers/dates), while leaving the cell style unchanged.
function Set-TextValue {
    param($range, $value)
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '26.854.44'
$ws.Range('E2').Value = '  -1.69%  '
# Row 3
Set-TextValue 'D3' '1.828.82'
$ws.Range('E3').Value = '  -1.36%  '
# Row 4
$ws.Range('E4').Value = '  +0.67%  '
# Row 5
Set-TextValue 'D5' '310.88'
$ws.Range('E5').Value = '  -1.04%  '
# Row 6
$ws.Range('E6').Value = '  +0.56%  '
# Row 7
Set-TextValue 'D7' '0.4578'
$ws.Range('E7').Value = '  -0.70%  '
# Row 8
Set-TextValue 'D8' '0.3674'
$ws.Range('E8').Value = '  -0.99%  '
# Row 9
Set-TextValue 'D9' '0.07167'
$ws.Range('E9').Value = '  -1.89%  '
# Row 10
Set-TextValue 'D10' '0.8742'
$ws.Range('E10').Value = '  -1.33%  '
# Row 11
Set-TextValue 'D11' '0.07817'
$ws.Range('E11').Value = '  +0.39%  '
# Row 12
Set-TextValue 'D12' '19.54'
$ws.Range('E12').Value = '  -2.01%  '
# Row 13
Set-TextValue 'D13' '1.869.38'
$ws.Range('E13').Value = '  +1.28%  '
# Row 14
Set-TextValue 'D14' '5.321'
$ws.Range('E14').Value = '  -1.01%  '
# Row 15
Set-TextValue 'D15' '6.365'
$ws.Range('E15').Value = '  -2.84%  '
# Row 16
Set-TextValue 'D16' '87.09'
$ws.Range('E16').Value = '  -5.02%  '
# Row 17
$ws.Range('E17').Value = '  +0.65%  '
# Row 18
Set-TextValue 'D18' '0.000008710'
$ws.Range('E18').Value = '  -3.00%  '
# Row 19
$ws.Range('E19').Value = '  +0.54%  '
# Row 20
Set-TextValue 'D20' '26.880.00'
$ws.Range('E20').Value = '  -1.62%  '
# Row 21
Set-TextValue 'D21' '14.47'
$ws.Range('E21').Value = '  -2.10%  '
# Row 22
Set-TextValue 'D22' '4.985'
$ws.Range('E22').Value = '  -2.81%  '
# Row 23
Set-TextValue 'D23' '10.45'
$ws.Range('E23').Value = '  -0.56%  '
# Row 24
Set-TextValue 'D24' '1.996'
$ws.Range('E24').Value = '  +3.63%  '
# Row 25
Set-TextValue 'D25' '151.69'
$ws.Range('E25').Value = '  +0.10%  '
# Row 26
Set-TextValue 'D26' '18.18'
$ws.Range('E26').Value = '  -1.22%  '
# Row 27
Set-TextValue 'D27' '1.992'
$ws.Range('E27').Value = '  -3.35%  '
# Row 28
Set-TextValue 'D28' '113.57'
$ws.Range('E28').Value = '  -2.28%  '
# Row 29
Set-TextValue 'D29' '4.919'
# Row 30
Set-TextValue 'D30' '0.08801'
$ws.Range('E30').Value = '  -0.37%  '
# Row 31
Set-TextValue 'D31' '3.100'
$ws.Range('E31').Value = '  -0.95%  '
# Row 32
Set-TextValue 'D32' '0.7414'
$ws.Range('E32').Value = '  -4.54%  '
# Row 33
Set-TextValue 'D33' '4.479'
$ws.Range('E33').Value = '  -0.43%  '
# Row 34
Set-TextValue 'D34' '1.129'
$ws.Range('E34').Value = '  -3.56%  '
# Row 35
Set-TextValue 'D35' '2.512'
$ws.Range('E35').Value = '  -5.68%  '
# Row 36
Set-TextValue 'D36' '1.082'
$ws.Range('E36').Value = '  +0.49%  '
# Row 37
Set-TextValue 'D37' '0.01936'
$ws.Range('E37').Value = '  -1.25%  '
# Row 38
Set-TextValue 'D38' '0.05128'
$ws.Range('E38').Value = '  -2.05%  '
# Row 39
Set-TextValue 'D39' '2.906'
# Row 40
Set-TextValue 'D40' '6.916'
$ws.Range('E40').Value = '  -1.09%  '
# Row 41
Set-TextValue 'D41' '0.4961'
$ws.Range('E41').Value = '  -3.58%  '
# Row 42
$ws.Range('E42').Value = '  -2.55%  '
# Row 43
Set-TextValue 'D43' '8.270'
$ws.Range('E43').Value = '  -1.75%  '
# Row 44
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D44' '1.007'
$ws.Range('E44').Value = '  +0.64%  '
# Row 45
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D45' '0.4656'
$ws.Range('E45').Value = '  -3.27%  '
# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '10.13'
$ws.Range('E46').Value = '  -1.10%  '
# Row 47
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D47' '103.15'
$ws.Range('E47').Value = '  +0.57%  '
# Row 48
Set-TextValue 'D48' '1.601'
$ws.Range('E48').Value = '  -3.00%  '
# Row 49
Set-TextValue 'D49' '0.06064'
$ws.Range('E49').Value = '  -2.48%  '
# Row 50
Set-TextValue 'D50' '64.71'
$ws.Range('E50').Value = '  -1.31%  '
# Row 51
Set-TextValue 'D51' '36.56'
$ws.Range('E51').Value = '  -0.53%  '
